$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.506869
$ws.Range("H2").Value = 49.520607
$ws.Range("I2").Value = 0.2165594803671733
$ws.Range("J2").Value = 0.2165594803671733
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 68.00339
$ws.Range("N2").Value = 204.01017
$ws.Range("O2").Value = 0.6265962299909886
$ws.Range("P2").Value = 0.6265962299909885
$ws.Range("Q2").Value = 1122.52305028591
$ws.Range("R2").Value = 10102.70745257319
$ws.Range("S2").Value = 0.1356953539668783
$ws.Range("T2").Value = 0.1356953539668783

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.506869
$ws.Range("H3").Value = 49.520607
$ws.Range("I3").Value = 0.2165594803671733
$ws.Range("J3").Value = 0.2165594803671733
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.991529999999999
$ws.Range("N3").Value = 26.97459
$ws.Range("O3").Value = 0.08284967558015671
$ws.Range("P3").Value = 0.08284967558015671
$ws.Range("Q3").Value = 148.42200781957
$ws.Range("R3").Value = 1335.79807037613
$ws.Range("S3").Value = 0.01794188269222763
$ws.Range("T3").Value = 0.01794188269222763

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.506869
$ws.Range("H4").Value = 49.520607
$ws.Range("I4").Value = 0.2165594803671733
$ws.Range("J4").Value = 0.2165594803671733
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.76843933333333
$ws.Range("N4").Value = 32.305318
$ws.Range("O4").Value = 0.09922245772090688
$ws.Range("P4").Value = 0.09922245772090688
$ws.Range("Q4").Value = 177.7532174097807
$ws.Range("R4").Value = 1599.778956688026
$ws.Range("S4").Value = 0.02148756388479341
$ws.Range("T4").Value = 0.02148756388479342

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.506869
$ws.Range("H5").Value = 49.520607
$ws.Range("I5").Value = 0.2165594803671733
$ws.Range("J5").Value = 0.2165594803671733
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 20.764887
$ws.Range("N5").Value = 62.294661
$ws.Range("O5").Value = 0.1913316367079478
$ws.Range("P5").Value = 0.1913316367079478
$ws.Range("Q5").Value = 342.7632695088029
$ws.Range("R5").Value = 3084.869425579227
$ws.Range("S5").Value = 0.04143467982327395
$ws.Range("T5").Value = 0.04143467982327396

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 24.781512
$ws.Range("H6").Value = 74.34453600000001
$ws.Range("I6").Value = 0.3251174623990092
$ws.Range("J6").Value = 0.3251174623990092
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 68.00339
$ws.Range("N6").Value = 204.01017
$ws.Range("O6").Value = 0.6265962299909886
$ws.Range("P6").Value = 0.6265962299909885
$ws.Range("Q6").Value = 1685.22682532568
$ws.Range("R6").Value = 15167.04142793112
$ws.Range("S6").Value = 0.2037173762434561
$ws.Range("T6").Value = 0.2037173762434561

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 24.781512
$ws.Range("H7").Value = 74.34453600000001
$ws.Range("I7").Value = 0.3251174623990092
$ws.Range("J7").Value = 0.3251174623990092
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.991529999999999
$ws.Range("N7").Value = 26.97459
$ws.Range("O7").Value = 0.08284967558015671
$ws.Range("P7").Value = 0.08284967558015671
$ws.Range("Q7").Value = 222.82370859336
$ws.Range("R7").Value = 2005.41337734024
$ws.Range("S7").Value = 0.02693587628520171
$ws.Range("T7").Value = 0.02693587628520171

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 24.781512
$ws.Range("H8").Value = 74.34453600000001
$ws.Range("I8").Value = 0.3251174623990092
$ws.Range("J8").Value = 0.3251174623990092
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.76843933333333
$ws.Range("N8").Value = 32.305318
$ws.Range("O8").Value = 0.09922245772090688
$ws.Range("P8").Value = 0.09922245772090688
$ws.Range("Q8").Value = 266.858208560272
$ws.Range("R8").Value = 2401.723877042448
$ws.Range("S8").Value = 0.03225895366721422
$ws.Range("T8").Value = 0.03225895366721422

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 24.781512
$ws.Range("H9").Value = 74.34453600000001
$ws.Range("I9").Value = 0.3251174623990092
$ws.Range("J9").Value = 0.3251174623990092
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.764887
$ws.Range("N9").Value = 62.294661
$ws.Range("O9").Value = 0.1913316367079478
$ws.Range("P9").Value = 0.1913316367079478
$ws.Range("Q9").Value = 514.585296369144
$ws.Range("R9").Value = 4631.267667322296
$ws.Range("S9").Value = 0.0622052562031371
$ws.Range("T9").Value = 0.0622052562031371

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 24.67943933333333
$ws.Range("H10").Value = 74.038318
$ws.Range("I10").Value = 0.3237783348120013
$ws.Range("J10").Value = 0.3237783348120013
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 68.00339
$ws.Range("N10").Value = 204.01017
$ws.Range("O10").Value = 0.6265962299909886
$ws.Range("P10").Value = 0.6265962299909885
$ws.Range("Q10").Value = 1678.285537966007
$ws.Range("R10").Value = 15104.56984169406
$ws.Range("S10").Value = 0.2028782839459601
$ws.Range("T10").Value = 0.2028782839459601

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 24.67943933333333
$ws.Range("H11").Value = 74.038318
$ws.Range("I11").Value = 0.3237783348120013
$ws.Range("J11").Value = 0.3237783348120013
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.991529999999999
$ws.Range("N11").Value = 26.97459
$ws.Range("O11").Value = 0.08284967558015671
$ws.Range("P11").Value = 0.08284967558015671
$ws.Range("Q11").Value = 221.9059191488467
$ws.Range("R11").Value = 1997.15327233962
$ws.Range("S11").Value = 0.02682492999905767
$ws.Range("T11").Value = 0.02682492999905767

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 24.67943933333333
$ws.Range("H12").Value = 74.038318
$ws.Range("I12").Value = 0.3237783348120013
$ws.Range("J12").Value = 0.3237783348120013
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.76843933333333
$ws.Range("N12").Value = 32.305318
$ws.Range("O12").Value = 0.09922245772090688
$ws.Range("P12").Value = 0.09922245772090688
$ws.Range("Q12").Value = 265.7590452416804
$ws.Range("R12").Value = 2391.831407175124
$ws.Range("S12").Value = 0.03212608213682943
$ws.Range("T12").Value = 0.03212608213682943

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 24.67943933333333
$ws.Range("H13").Value = 74.038318
$ws.Range("I13").Value = 0.3237783348120013
$ws.Range("J13").Value = 0.3237783348120013
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 20.764887
$ws.Range("N13").Value = 62.294661
$ws.Range("O13").Value = 0.1913316367079478
$ws.Range("P13").Value = 0.1913316367079478
$ws.Range("Q13").Value = 512.4657689800219
$ws.Range("R13").Value = 4612.191920820198
$ws.Range("S13").Value = 0.06194903873015412
$ws.Range("T13").Value = 0.06194903873015412

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.255437
$ws.Range("H14").Value = 30.766311
$ws.Range("I14").Value = 0.1345447224218162
$ws.Range("J14").Value = 0.1345447224218162
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 68.00339
$ws.Range("N14").Value = 204.01017
$ws.Range("O14").Value = 0.6265962299909886
$ws.Range("P14").Value = 0.6265962299909885
$ws.Range("Q14").Value = 697.40448193143
$ws.Range("R14").Value = 6276.64033738287
$ws.Range("S14").Value = 0.08430521583469407
$ws.Range("T14").Value = 0.08430521583469407

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.255437
$ws.Range("H15").Value = 30.766311
$ws.Range("I15").Value = 0.1345447224218162
$ws.Range("J15").Value = 0.1345447224218162
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 8.991529999999999
$ws.Range("N15").Value = 26.97459
$ws.Range("O15").Value = 0.08284967558015671
$ws.Range("P15").Value = 0.08284967558015671
$ws.Range("Q15").Value = 92.21206944860999
$ws.Range("R15").Value = 829.90862503749
$ws.Range("S15").Value = 0.01114698660366971
$ws.Range("T15").Value = 0.01114698660366971

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.255437
$ws.Range("H16").Value = 30.766311
$ws.Range("I16").Value = 0.1345447224218162
$ws.Range("J16").Value = 0.1345447224218162
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.76843933333333
$ws.Range("N16").Value = 32.305318
$ws.Range("O16").Value = 0.09922245772090688
$ws.Range("P16").Value = 0.09922245772090688
$ws.Range("Q16").Value = 110.435051171322
$ws.Range("R16").Value = 993.915460541898
$ws.Range("S16").Value = 0.01334985803206981
$ws.Range("T16").Value = 0.01334985803206981

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.255437
$ws.Range("H17").Value = 30.766311
$ws.Range("I17").Value = 0.1345447224218162
$ws.Range("J17").Value = 0.1345447224218162
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 20.764887
$ws.Range("N17").Value = 62.294661
$ws.Range("O17").Value = 0.1913316367079478
$ws.Range("P17").Value = 0.1913316367079478
$ws.Range("Q17").Value = 212.952990440619
$ws.Range("R17").Value = 1916.576913965571
$ws.Range("S17").Value = 0.02574266195138261
$ws.Range("T17").Value = 0.02574266195138262
